$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at row 3 (shifts old rows 3-12 down to 4-13) ---
$ws.Rows.Item(3).Insert()

# --- 2. Populate the new row 3 ("No Tax Address" record) ---
$ws.Range("A3").Value = "No Tax Address"
$ws.Range("F3").Value = "Test"
$ws.Range("G3").Value = "Qa"
$ws.Range("L3").Value = "844 N Colony Rd"
$ws.Range("M3").Value = "Wallingford"
$ws.Range("N3").Value = "United States"
$ws.Range("O3").Value = "Connecticut"
$ws.Range("O3").Style = "Normal"
$ws.Range("P3").Value = "'06492"
$ws.Range("Q3").Value = 9898989898

# --- 3. Replace the old "2 Dream Valley Drive" address with the new
#        "844 North Court" address everywhere it is used (row 2, and
#        what is now row 13 after the insert above) ---
$ws.Range("L2").Value = "844 North Court"
$ws.Range("M2").Value = "Albany"
$ws.Range("O2").Value = "New York"
$ws.Range("O2").Style = "Normal"
$ws.Range("P2").Value = "'12211"

$ws.Range("L13").Value = "844 North Court"
$ws.Range("M13").Value = "Albany"
$ws.Range("O13").Value = "New York"
$ws.Range("O13").Style = "Normal"
$ws.Range("P13").Value = "'12211"

# --- 4. Rebuild the hyperlinks collection with the shifted addresses.
#        (Row insert does not relocate existing hyperlink anchors, so
#        the three below K9/B11/D11 must move to K10/B12/D12.) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Lotuswave@123")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:avayugundla@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:avayugundla@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Lotuswave@123")
$ws.Hyperlinks.Add($ws.Range("K10"), "mailto:vnarra@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:avayugundla@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:Lotuswave@123")

# --- 5. Update the sheet view: scroll so column E is the leftmost
#        visible column, and select the newly-relevant L13:P13 block ---
$ws.Range("L13:P13").Select()
$excel.ActiveWindow.ScrollColumn = 5
